# Add "Booster Sample No" and "Is Standard Rotation" columns with data,
# set sample sizes and whether each set is in the current standard rotation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cells (bold font, thin border, matching the existing header style) ---
$ws.Range("C1").Value = "Booster Sample No"
$ws.Range("D1").Value = "Is Standard Rotation"

$headerRange = $ws.Range("C1:D1")
$headerRange.Font.Bold = $true
$headerRange.Borders.LineStyle = 1

# --- Data: booster-pack sample size (numeric where known, "NNNN+" text otherwise) ---
$sampleNo = @(
    7020, 2088, 2010, 1536, 1530, 1674, 1264, 1480, 1274, 872,
    3888, 3672, 1580, 3136, 4628, 2736, 5040, 3192, 2187,
    "8000+", "8000+", "8000+", "8000+", "1900+", "8000+", "8000+", "8000+",
    "1500+", "8000+", "1500+", "8000+", "8000+", "8000+"
)

# --- Data: is the set currently part of the Standard rotation? ---
$isStandardRotation = @(
    0,0,0,0,0,0,0,0,0,0,
    0,0,0,0,0,0,0,0,0,
    0,1,1,1,1,1,1,1,1,1,1,1,1,1
)

for ($i = 0; $i -lt $sampleNo.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $sampleNo[$i]
    $ws.Cells.Item($row, 4).Value = $isStandardRotation[$i]
}

# Gray-fill style carried by the numeric "sample size" cells (rows 2-21)
$ws.Range("C2:C21").Interior.Pattern = 0

# --- Column widths ---
$ws.Columns.Item(2).ColumnWidth = 10.08984375
$ws.Columns.Item(3).ColumnWidth = 16.81640625
$ws.Columns.Item(4).ColumnWidth = 18.36328125

# --- View: zoom + selection to match the edited workbook ---
$ws.Application.ActiveWindow.Zoom = 85
$ws.Range("C2:C21").Select()
